$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1) - update F column ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13820
$ws1.Range("F8").Value = 1032
$ws1.Range("F9").Value = 13899
$ws1.Range("F10").Value = 14777
$ws1.Range("F20").Value = 20
$ws1.Range("F23").Value = 1149
$ws1.Range("F26").Value = 5738
$ws1.Range("F27").Value = 945
$ws1.Range("F28").Value = 1057
$ws1.Range("F29").Value = 5416
$ws1.Range("F31").Value = 50
$ws1.Range("F32").Value = 270

# "全部类型" sheet (sheet4) - same underlying rows, offset by one extra row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13820
$ws4.Range("F9").Value = 1032
$ws4.Range("F10").Value = 13899
$ws4.Range("F11").Value = 14777
$ws4.Range("F21").Value = 20
$ws4.Range("F24").Value = 1149
$ws4.Range("F27").Value = 5738
$ws4.Range("F28").Value = 945
$ws4.Range("F29").Value = 1057
$ws4.Range("F30").Value = 5416
$ws4.Range("F32").Value = 50
$ws4.Range("F33").Value = 270
